# entry 6 compound equivalents updated along with new prediction results
#
# Row 6 ("Benzoic acid" + "Benzenamine hydrochloride (1:1)") gets its amine
# equivalents, activator equivalents and base equivalents revised, the
# amine name corrected to drop the stray "(1:1)" suffix, and the model's
# predicted reaction yield refreshed to match the new inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Amine (equiv.) B6: 1 -> 1.6
$ws.Range("B6").Value = 1.6

# Amine (name) C6: "Benzenamine hydrochloride (1:1)" -> "Benzenamine hydrochloride"
$ws.Range("C6").Value = "Benzenamine hydrochloride"

# Activator (equiv.) F6: 1.5 -> 2.5
$ws.Range("F6").Value = 2.5

# Base (equiv.) H6: 2 -> 3.3
$ws.Range("H6").Value = 3.3

# Predicted Reaction_Yield N6: 89.904756000000006 -> 72.972087999999999
$ws.Range("N6").Value = 72.972087999999999

# Scroll the view over and move the selection to N9, matching where the
# author was working when the prediction results were refreshed.
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N9").Select()
